$d = $word.ActiveDocument

# Locate the paragraph that currently holds "This is the left laptop " plus the
# _GoBack bookmark, and insert 11 blank paragraphs plus a new paragraph with
# the "This is testing 2 by right laptop" text right after it. The bookmark
# ends up attached to the newly added paragraph (it travels with the
# following content since it is inserted right after the existing run).

$target = $d.Paragraphs(3).Range

$p = $target.Paragraphs(1).Range
$p.InsertParagraphAfter()

# Re-fetch the paragraph collection; paragraph 3's range now ends just after
# the first blank paragraph break. Insert 10 more blank paragraphs.
for ($i = 0; $i -lt 10; $i++) {
    $d.Paragraphs(4).Range.InsertParagraphAfter()
}

# Finally insert the new text paragraph.
$d.Paragraphs(14).Range.InsertParagraphAfter()
$d.Paragraphs(15).Range.Text = "This is testing 2 by right laptop"
